# refactor tests to that tests can be ran consecutively
# Update mock student data in StudentData.xlsx so that re-running the
# tests (which likely mutate/consume student records) doesn't collide
# with data left over from a previous run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: student S12345678A/Balqis -> S12345670A/Student 1
$ws.Range("A2").Value = "S12345670A"
$ws.Range("B2").Value = "Student 1"

# Row 3: student S12345678B/Lala -> S12345671B (name unchanged)
$ws.Range("A3").Value = "S12345671B"

# Row 4: student S12345678C/Pooh -> S12345677H/Student 8
$ws.Range("A4").Value = "S12345677H"
$ws.Range("B4").Value = "Student 8"

# Row 5 (S12345678D / Kinanti / SEG / Unassigned) is unchanged.

$wb.Save()
